$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value2 = "27.632.55"
$ws.Range("E2").Value2 = "  -0.71%  "
$ws.Range("D3").Value2 = "1.867.23"
$ws.Range("E3").Value2 = "  -1.04%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value2 = "1.014"
$ws.Range("E4").Value2 = "  -0.29%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value2 = "335.49"
$ws.Range("E5").Value2 = "  +0.24%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value2 = "1.013"
$ws.Range("E6").Value2 = "  -0.31%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value2 = "0.4665"
$ws.Range("E7").Value2 = "  -0.36%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value2 = "0.3926"
$ws.Range("E8").Value2 = "  +0.29%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value2 = "46.05"
$ws.Range("E9").Value2 = "  -2.91%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value2 = "0.07988"
$ws.Range("E10").Value2 = "  -0.85%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value2 = "1.002"
$ws.Range("E11").Value2 = "  -1.38%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value2 = "21.73"
$ws.Range("E12").Value2 = "  -1.02%  "
$ws.Range("D13").Value2 = "1.867.97"
$ws.Range("E13").Value2 = "  -2.47%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value2 = "5.955"
$ws.Range("E14").Value2 = "  -0.06%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value2 = "7.226"
$ws.Range("E15").Value2 = "  +2.01%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value2 = "1.015"
$ws.Range("E16").Value2 = "  -0.35%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value2 = "88.58"
$ws.Range("E17").Value2 = "  +1.51%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value2 = "0.06711"
$ws.Range("E18").Value2 = "  -0.82%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value2 = "0.00001045"
$ws.Range("E19").Value2 = "  -0.36%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value2 = "17.24"
$ws.Range("E20").Value2 = "  +0.08%  "
$ws.Range("E21").Value2 = "  -0.15%  "
$ws.Range("D22").Value2 = "27.637.71"
$ws.Range("E23").Value2 = "  -0.65%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value2 = "10.94"
$ws.Range("E24").Value2 = "  -0.03%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value2 = "2.312"
$ws.Range("E25").Value2 = "  -1.31%  "
$ws.Range("D26").Value2 = "2.091.06"
$ws.Range("E26").Value2 = "  -2.22%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value2 = "159.09"
$ws.Range("E27").Value2 = "  -0.51%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value2 = "19.67"
$ws.Range("E28").Value2 = "  -1.96%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value2 = "2.138"
$ws.Range("E29").Value2 = "  +2.56%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value2 = "5.439"
$ws.Range("E30").Value2 = "  -0.54%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value2 = "121.67"
$ws.Range("E31").Value2 = "  -0.14%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value2 = "0.9773"
$ws.Range("E32").Value2 = "  +0.83%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value2 = "0.09461"
$ws.Range("E33").Value2 = "  -0.22%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value2 = "3.623"
$ws.Range("E34").Value2 = "  -0.63%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value2 = "5.313"
$ws.Range("E35").Value2 = "  -0.74%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value2 = "1.341"
$ws.Range("E36").Value2 = "  -5.23%  "
$ws.Range("B37").Value2 = "VeChain"
$ws.Range("C37").Value2 = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value2 = "0.02237"
$ws.Range("E37").Value2 = "  -0.83%  "
$ws.Range("B38").Value2 = "Hedera"
$ws.Range("C38").Value2 = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value2 = "0.06031"
$ws.Range("E38").Value2 = "  -1.46%  "
$ws.Range("B39").Value2 = "FraxShare"
$ws.Range("C39").Value2 = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value2 = "8.345"
$ws.Range("E39").Value2 = "  +4.24%  "
$ws.Range("B40").Value2 = "TrustWalletToken"
$ws.Range("C40").Value2 = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value2 = "1.197"
$ws.Range("E40").Value2 = "  -1.31%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value2 = "1.012"
$ws.Range("E41").Value2 = "  -0.37%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value2 = "0.5954"
$ws.Range("E42").Value2 = "  -0.40%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value2 = "0.1874"
$ws.Range("E43").Value2 = "  -0.62%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value2 = "10.33"
$ws.Range("E44").Value2 = "  +0.55%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value2 = "1.250"
$ws.Range("E45").Value2 = "  -1.22%  "
$ws.Range("E46").Value2 = "  -0.95%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value2 = "12.15"
$ws.Range("E47").Value2 = "  -0.22%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value2 = "1.926"
$ws.Range("E48").Value2 = "  +0.04%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value2 = "0.06744"
$ws.Range("E49").Value2 = "  -2.53%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value2 = "111.83"
$ws.Range("E50").Value2 = "  -1.73%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value2 = "1.052"
$ws.Range("E51").Value2 = "  -1.60%  "
